$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.885.05'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '3.783.31'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.37'
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '163.57'
$ws.Range('E6').Value = '  -1.53%  '
$ws.Range('D7').Value = '3.780.36'
$ws.Range('E7').Value = '  -0.91%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('E10').Value = '  -1.91%  '
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.84'
$ws.Range('E12').Value = '  +8.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000247'
$ws.Range('E13').Value = '  -2.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.04'
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('D15').Value = '4.418.69'
$ws.Range('E15').Value = '  -0.87%  '
$ws.Range('D16').Value = '3.788.76'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').Value = '67.862.88'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.17'
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('E19').Value = '  +1.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.99'
$ws.Range('E20').Value = '  -1.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '458.57'
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.45'
$ws.Range('E22').Value = '  -4.38%  '
$ws.Range('E23').Value = '  -1.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.26'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  -2.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.88'
$ws.Range('E26').Value = '  -1.52%  '
$ws.Range('E27').Value = '  -1.10%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.90'
$ws.Range('E29').Value = '  -1.48%  '
$ws.Range('D30').Value = '3.933.36'
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('E31').Value = '  -6.85%  '
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('E33').Value = '  -1.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.98'
$ws.Range('E34').Value = '  -1.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '8.92'
$ws.Range('E36').Value = '  -1.75%  '
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.147'
$ws.Range('E38').Value = '  +7.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.80'
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('E40').Value = '  -2.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.978'
$ws.Range('E41').Value = '  -1.90%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.56'
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.07'
$ws.Range('E45').Value = '  -2.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '152.50'
$ws.Range('E46').Value = '  +2.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.294'
$ws.Range('E47').Value = '  -1.94%  '
$ws.Range('E48').Value = '  -1.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.30'
$ws.Range('E49').Value = '  -0.31%  '
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.40'
$ws.Range('E51').Value = '  -7.45%  '
